$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60 (shifts existing rows 60-94 down to 61-95)
$ws.Rows.Item(60).Insert()

# Populate the newly inserted row 60 with its data
$ws.Cells.Item(60, 1).Value = 1
$ws.Cells.Item(60, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(60, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(60, 4).Value = 45097
$ws.Cells.Item(60, 5).Value = 15
$ws.Cells.Item(60, 6).Value = 100112009
$ws.Cells.Item(60, 7).Value = "Acelga"
$ws.Cells.Item(60, 8).Value = "Sin especificar"
$ws.Cells.Item(60, 9).Value = "Segunda"
$ws.Cells.Item(60, 10).Value = 300
$ws.Cells.Item(60, 11).Value = 1000
$ws.Cells.Item(60, 12).Value = 1500
$ws.Cells.Item(60, 13).Value = 1250
$ws.Cells.Item(60, 14).Value = "`$/atado 2,5 a 3 kilos"
$ws.Cells.Item(60, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(60, 16).Value = 417
$ws.Cells.Item(60, 17).Value = 3
$ws.Cells.Item(60, 18).Value = "Hortaliza"
